$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1 "Tipo", matching the style of the other header cells (A1:C1)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Tipo"

# Update existing prediction values in B2 and C2
$ws.Range("B2").Value = 0.2383084938879716
$ws.Range("C2").Value = 0.9953422544939506

# Add new data cell D2 "single"
$ws.Range("D2").Value = "single"
